$d = $word.ActiveDocument

$replacements = @(
    @{Old = "46÷9="; New = "56÷9="},
    @{Old = "19÷5="; New = "20÷4="},
    @{Old = "50÷3="; New = "81÷8="},
    @{Old = "97÷6="; New = "46÷5="},
    @{Old = "61÷6="; New = "47÷5="},
    @{Old = "14÷5="; New = "18÷2="},
    @{Old = "51÷2="; New = "35÷2="},
    @{Old = "34÷2="; New = "32÷2="},
    @{Old = "81÷6="; New = "31÷8="},
    @{Old = "56÷2="; New = "85÷5="},
    @{Old = "39÷4="; New = "95÷6="},
    @{Old = "50÷2="; New = "48÷6="},
    @{Old = "11÷9="; New = "63÷9="},
    @{Old = "39÷6="; New = "97÷2="},
    @{Old = "55÷3="; New = "52÷3="},
    @{Old = "89÷4="; New = "13÷2="},
    @{Old = "40÷3="; New = "95÷7="},
    @{Old = "54÷5="; New = "48÷7="},
    @{Old = "38÷5="; New = "86÷6="},
    @{Old = "34÷3="; New = "60÷6="},
    @{Old = "49÷3="; New = "40÷8="},
    @{Old = "29÷4="; New = "16÷9="},
    @{Old = "14÷9="; New = "18÷4="},
    @{Old = "51÷7="; New = "57÷3="},
    @{Old = "85÷8="; New = "22÷8="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $true, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
